# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 48, shifting the existing rows
# 48-105 down to 49-106.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(48).Insert()

$ws.Cells.Item(48, 1).Value = 11
$ws.Cells.Item(48, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(48, 3).Value = "Bíobío"
$ws.Cells.Item(48, 4).Value = 44483
$ws.Cells.Item(48, 5).Value = 8
$ws.Cells.Item(48, 6).Value = "Fruta"
$ws.Cells.Item(48, 7).Value = 100108
$ws.Cells.Item(48, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(48, 9).Value = 100108005
$ws.Cells.Item(48, 10).Value = "Piña"
$ws.Cells.Item(48, 11).Value = "Caramelo"
$ws.Cells.Item(48, 12).Value = "Segunda"
$ws.Cells.Item(48, 13).Value = 200
$ws.Cells.Item(48, 14).Value = 21000
$ws.Cells.Item(48, 15).Value = 22000
$ws.Cells.Item(48, 16).Value = 21500
$ws.Cells.Item(48, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(48, 18).Value = "Ecuador"
$ws.Cells.Item(48, 19).Value = 1536
$ws.Cells.Item(48, 20).Value = 14
